# Updated cryptos list on Sat Mar 11 09:43:39 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 40/41 (Aptos / Frax) swap places along with updated values.
# Leading "'" forces text storage so numeric-looking prices are not
# reinterpreted as numbers by Excel.
$ws.Cells.Item(40, 2).Value = "Frax"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(40, 4).Value = "'0.9354"
$ws.Cells.Item(40, 5).Value = "  -6.56%  "

$ws.Cells.Item(41, 2).Value = "Aptos"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(41, 4).Value = "'10.14"
$ws.Cells.Item(41, 5).Value = "  +0.56%  "

# Rows 44/45 (TheSandbox / PancakeSwap) swap places along with updated values.
$ws.Cells.Item(44, 2).Value = "PancakeSwap"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(44, 4).Value = "'3.514"
$ws.Cells.Item(44, 5).Value = "  +0.73%  "

$ws.Cells.Item(45, 2).Value = "TheSandbox"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(45, 4).Value = "'0.5230"
$ws.Cells.Item(45, 5).Value = "  +0.19%  "

# Price / Volume(1h) updates for all other rows.
$ws.Cells.Item(2, 4).Value = "'20.153.82"
$ws.Cells.Item(2, 5).Value = "  +1.24%  "
$ws.Cells.Item(3, 4).Value = "'1.441.23"
$ws.Cells.Item(3, 5).Value = "  +2.69%  "
$ws.Cells.Item(4, 4).Value = "'1.009"
$ws.Cells.Item(4, 5).Value = "  +0.65%  "
$ws.Cells.Item(5, 4).Value = "'0.9118"
$ws.Cells.Item(5, 5).Value = "  -8.98%  "
$ws.Cells.Item(6, 4).Value = "'276.97"
$ws.Cells.Item(6, 5).Value = "  +2.70%  "
$ws.Cells.Item(7, 4).Value = "'0.3658"
$ws.Cells.Item(7, 5).Value = "  -0.29%  "
$ws.Cells.Item(8, 4).Value = "'0.3126"
$ws.Cells.Item(8, 5).Value = "  +2.90%  "
$ws.Cells.Item(9, 4).Value = "'38.78"
$ws.Cells.Item(9, 5).Value = "  -0.82%  "
$ws.Cells.Item(10, 4).Value = "'1.018"
$ws.Cells.Item(10, 5).Value = "  +4.85%  "
$ws.Cells.Item(11, 4).Value = "'0.06522"
$ws.Cells.Item(11, 5).Value = "  +1.04%  "
$ws.Cells.Item(12, 4).Value = "'1.004"
$ws.Cells.Item(12, 5).Value = "  +0.06%  "
$ws.Cells.Item(13, 4).Value = "'5.387"
$ws.Cells.Item(13, 5).Value = "  +2.26%  "
$ws.Cells.Item(14, 4).Value = "'17.50"
$ws.Cells.Item(14, 5).Value = "  +5.56%  "
$ws.Cells.Item(15, 4).Value = "'6.063"
$ws.Cells.Item(15, 5).Value = "  +0.15%  "
$ws.Cells.Item(16, 4).Value = "'1.444.15"
$ws.Cells.Item(16, 5).Value = "  +2.41%  "
$ws.Cells.Item(17, 4).Value = "'0.00001014"
$ws.Cells.Item(17, 5).Value = "  +1.42%  "
$ws.Cells.Item(18, 4).Value = "'0.9382"
$ws.Cells.Item(18, 5).Value = "  -6.33%  "
$ws.Cells.Item(19, 4).Value = "'0.05638"
$ws.Cells.Item(19, 5).Value = "  -0.56%  "
$ws.Cells.Item(20, 4).Value = "'67.70"
$ws.Cells.Item(20, 5).Value = "  -5.51%  "
$ws.Cells.Item(21, 4).Value = "'5.409"
$ws.Cells.Item(21, 5).Value = "  -1.40%  "
$ws.Cells.Item(22, 4).Value = "'14.45"
$ws.Cells.Item(22, 5).Value = "  +2.24%  "
$ws.Cells.Item(23, 4).Value = "'10.80"
$ws.Cells.Item(23, 5).Value = "  +2.16%  "
$ws.Cells.Item(24, 4).Value = "'2.267"
$ws.Cells.Item(24, 5).Value = "  +0.00%  "
$ws.Cells.Item(25, 4).Value = "'20.183.15"
$ws.Cells.Item(25, 5).Value = "  +1.27%  "
$ws.Cells.Item(26, 4).Value = "'2.182"
$ws.Cells.Item(26, 5).Value = "  -0.76%  "
$ws.Cells.Item(27, 4).Value = "'136.25"
$ws.Cells.Item(27, 5).Value = "  +0.31%  "
$ws.Cells.Item(28, 4).Value = "'16.97"
$ws.Cells.Item(28, 5).Value = "  +2.77%  "
$ws.Cells.Item(29, 4).Value = "'1.598.81"
$ws.Cells.Item(29, 5).Value = "  +2.04%  "
$ws.Cells.Item(30, 4).Value = "'110.81"
$ws.Cells.Item(30, 5).Value = "  +3.64%  "
$ws.Cells.Item(31, 4).Value = "'3.738"
$ws.Cells.Item(31, 5).Value = "  -3.11%  "
$ws.Cells.Item(32, 4).Value = "'0.8058"
$ws.Cells.Item(32, 5).Value = "  +1.11%  "
$ws.Cells.Item(33, 4).Value = "'4.815"
$ws.Cells.Item(33, 5).Value = "  -7.14%  "
$ws.Cells.Item(34, 4).Value = "'0.07701"
$ws.Cells.Item(34, 5).Value = "  +0.98%  "
$ws.Cells.Item(35, 4).Value = "'0.05994"
$ws.Cells.Item(35, 5).Value = "  +4.39%  "
$ws.Cells.Item(36, 4).Value = "'1.465"
$ws.Cells.Item(36, 5).Value = "  +9.57%  "
$ws.Cells.Item(37, 4).Value = "'4.690"
$ws.Cells.Item(37, 5).Value = "  +0.08%  "
$ws.Cells.Item(38, 4).Value = "'1.133"
$ws.Cells.Item(38, 5).Value = "  +7.76%  "
$ws.Cells.Item(39, 4).Value = "'0.01992"
$ws.Cells.Item(39, 5).Value = "  -0.80%  "
$ws.Cells.Item(42, 4).Value = "'0.1832"
$ws.Cells.Item(42, 5).Value = "  -3.97%  "
$ws.Cells.Item(43, 4).Value = "'7.189"
$ws.Cells.Item(43, 5).Value = "  -13.39%  "
$ws.Cells.Item(46, 4).Value = "'11.99"
$ws.Cells.Item(46, 5).Value = "  +0.76%  "
$ws.Cells.Item(47, 4).Value = "'119.04"
$ws.Cells.Item(47, 5).Value = "  +9.19%  "
$ws.Cells.Item(48, 4).Value = "'0.5126"
$ws.Cells.Item(48, 5).Value = "  +2.10%  "
$ws.Cells.Item(49, 4).Value = "'1.765"
$ws.Cells.Item(49, 5).Value = "  +0.68%  "
$ws.Cells.Item(50, 4).Value = "'0.06322"
$ws.Cells.Item(50, 5).Value = "  +3.15%  "
$ws.Cells.Item(51, 4).Value = "'0.9959"
$ws.Cells.Item(51, 5).Value = "  -0.53%  "
